$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCRUM")

$ws.Range("H15").Value = "1 - Soucis Git et EC2 (toujours et encore)`n2 - EC2 - Finalement call prévu avec Greg demain`n3 - Ai regardé RunPod.io (alternative à EC2?). Pas tout compris à part qu'il fallait que je paie encore. Call Adrien demain ? Standby pour l'instant.`nPartie 2/n`n4 - Concernant les models... Heu je veux bien aider et rentrer dans la danse mais j'arrive avec 9 jours de retard... J'ai lu un peu. Ai voulu jeter un oeil sur le code de transfer learning de Patrice mais je l'ai reçu un peu tard sur GitHub (23H30). En plus je comprends pas tout ce qui est fait... Je me demande si demain après midi ca vaudrait pas le coup que Patrice nous présente (à toi et moi par exemple) son code de transfer learning (.py ou .ipynb on s'en fout on s'adaptera toujours). Faudra voir demain. By the way, je me demande si on pourrait pas essayer de faire du transfer learning avec un autre modèle histoire de varier les plaisirs : EfficientNetB7 ou ResNet101V2 par exemple. Pour l'instant dans le code de Patrice c'est du InceptionV3 je crois. A discuter. Encore une fois, je découvre le sujet et on est mercredi... En tout cas si je lance des runs les résultats et artifacts iront sur le serveur mlflow tracking dès le départ.`n5 - Ai regardé comment ajouter des images dans les artifacts de mlflow tracking (matrice de confusion par exemple)`n6 - Ai finalisé le slide archi (mais pas partagé sur github because j'ai peur de tout casser pour la 53eme fois)`n7 - Call avec Patrice qui avait des questions mlflow tracking (support jupyter notebook)`n8 - J'ai pas l'occasion de reparler des prédiction avec Quentin. J'imagine que tout tourne de son côté.`n9 - En attendant le code de Patrice, je suis revenu ce soir sur un bout de code flask qui date de ce WE : on choisit une image, elle s'affiche, on peut zoomer et  sélectionner une zone d'interêt avant de ne garder qu'une image 512x512.`n10 - Un modèle de transfer learning fonctionne`nAccuracy: 0.0694`nPrecision: 0.0634`nRecall: 0.0694`nF1 Score: 0.0576"
$ws.Range("J15").Value = "Voir sensibilité = capacité du test à bien détecter les gens vraiment malades.`nOn peut pas utiliser la sensibilité toute seule.`nRegarder aussi specificté`nLa sensibilité correspond à la probabilité d'être testé positif lorsqu'on est malade : La spécificité correspond à la probabilité d'être testé négatif lorsqu'on n'est pas malade"
$ws.Range("E15").Value = "Fine tuning`nJe cherche à améliorer la loss et accuracy`nAccuracy = 0.2 puis = 0.5`nLoss = ??`n"
$ws.Range("C15").Value = "Aurélie indique`nEtape 1 : Accuracy & F1 score*`nEtape 2 : Ensuite les autres metrique"
$ws.Range("G15").Value = "Regarder F1 score"
$ws.Range("F15").Value = "J'arrive pas à logguer les matrices de confusion`nJ'overfit pas mal`n"
$ws.Range("B15").Value = "Je travaille sur moins de classes`n3 classes à la fin`ncouche de convolution en plus`nLa loss baisse puis Plateau `n25 epochs`nJe commence à overfitter`nJe peux encore baisser (15 epochs)`nJe suis à 20 min 3 classes`nValidation`nAccuracy = 0,6 epoch 8`nLoss = 0,8 epoch 8"
$ws.Range("H16").Value = "1 - Le modèle InceptionV3 (brut de fonderie) envoit ses courbes, matrices, paramètres, résultats et benchs sur mlflow tracking`n2 - Slides relus et envoyés avant minuit"
$ws.Range("I16").Value = "EC2 & Git"
$ws.Range("J16").Value = "Revoir les slides et l'histoire. Le faire passer sur Google (attention animations). Prévoir un pdf ?`nPoint technos avec la dream team`nEntre temps, faire tourner le modèle pour diminuer l'overfit"
$ws.Range("H17").Value = "Pas de réunion SCRUM hier`nRépétiion en fin de journée `nModification du modèle`nAjout d'une fonction Early_Stopping`nA la sortie du InceptionV3`nSupprimer la couche GlobalAveragePooling2D et ajout d'une couche Flaten`nAjout d'une regulation L2 => Test associés`nAjout de 2 couches dense 1024 et 512 => Test associés`n"
$ws.Range("I17").Value = "EC2 & Git"
$ws.Range("J17").Value = "Répétition à 9H40`nConvertir les slides`nBien comprendre la nature des paramètres accuracy etc. qui sont sortis. C'est du train, valiadtion, pourquoi c'est si différent des graphes?`nContinuer à jourer avec le modèle`nDégeler des lots de 4 couches jusqu'à 10% du modèle InceptionV3 et voir ce que cela donne`n22 Mds de paramètres dans 311 couches actuellement`nChanger de modèle ?"

# Column width change: split col 9-10 group, widen column J (10) to 80
$ws.Columns.Item(10).ColumnWidth = 80

# Row heights follow Excel's autofit for wrapped text
$ws.Rows.Item(15).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(16).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(17).EntireRow.AutoFit() | Out-Null

# Update view: scroll position & selection on SCRUM sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("I17").Select()
